# "update after conceptual meeting"
# Simplify the `slug` (column B) values for the conference-output, thesis
# and report sub-types on the "nr" sheet: drop the now-redundant
# "conference-"/"-thesis"/"-report" segments since the parent-row code
# already conveys that context.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B15").Value = "theses"

$ws.Range("B11").Value = "paper"
$ws.Range("B12").Value = "proceeding"
$ws.Range("B13").Value = "programme"
$ws.Range("B14").Value = "poster"

$ws.Range("B16").Value = "bachelor"
$ws.Range("B17").Value = "master"
$ws.Range("B18").Value = "rigorous"
$ws.Range("B19").Value = "doctoral"
$ws.Range("B20").Value = "post-doctoral"

$ws.Range("B27").Value = "annual"
$ws.Range("B28").Value = "research"
$ws.Range("B29").Value = "project"
$ws.Range("B30").Value = "statistical-or-status"
$ws.Range("B31").Value = "conservation"
$ws.Range("B32").Value = "field"
$ws.Range("B33").Value = "business-trip"

# Leave the view with B33 selected (matches the cell last touched).
$ws.Range("B33").Select()
